# Fix the tab name "herbivoredata" -> "herbivoreData" and make that tab the
# active/selected tab (it had previously been left on "notes"), updating the
# frozen-pane selection on the herbivoreData sheet from F20 to F14.

$wb = $excel.ActiveWorkbook

# 1. Rename the mis-capitalized sheet.
$herbSheet = $wb.Worksheets.Item("herbivoredata")
$herbSheet.Name = "herbivoreData"

# 2. Activate it so it becomes the selected/visible tab (moves tabSelected
#    from "notes" to "herbivoreData" and updates the workbook's activeTab).
$herbSheet.Activate()

# 3. Update the active cell / selection within the frozen bottom-left pane.
$herbSheet.Range("F14").Select()
